# "Einnahme und Ausgaben korrrigiert"
# The "Zweifel Gutschrift Kiosk-Retouren" credit from Zweifel Pomy-Chips AG
# was booked as a negative expense on the "Ausgaben" sheet; it also needs
# to be recorded as a positive entry on the "Einnahmen" (income) table so
# both sheets reconcile correctly.

$wb = $excel.ActiveWorkbook
$wsAusgaben = $wb.Worksheets.Item("Ausgaben")
$ws = $wb.Worksheets.Item("Einnahmen")
$lo = $ws.ListObjects.Item(1)

# Give the income table a bit more visual emphasis (matches the table
# style used elsewhere for corrected / reconciled entries).
$lo.TableStyle = "TableStyleMedium9"

# Add the missing income row for the Kiosk credit note.
$newRow = $lo.ListRows.Add()
$r = $newRow.Range

$r.Item(1,1).Value = "Kiosk"
$r.Item(1,2).Value = "Zweifel Gutschrift Kiosk-Retouren"
$r.Item(1,3).Value = (Get-Date -Year 2023 -Month 12 -Day 8).Date
$r.Item(1,4).Value = 54.48
$r.Item(1,5).Value = "Zweifel Pomy-Chips AG"
$r.Item(1,6).Value = "Regensdorferstrasse 20, 8049 Zürich"

# Highlight the newly-added / corrected row so it's easy to spot.
$r.Interior.Pattern = 1
$r.Interior.ThemeColor = 5

$r.Item(1,3).NumberFormat = "dd/mm/yyyy;@"
$r.Item(1,4).NumberFormat = '"CHF" #,##0.00'

$topBorder = $r.Borders.Item(9)
$topBorder.LineStyle = 1
$topBorder.Weight = 2
$topBorder.ThemeColor = 5

$leftBorder = $r.Item(1,1).Borders.Item(7)
$leftBorder.LineStyle = 1
$leftBorder.Weight = 2
$leftBorder.ThemeColor = 5

# Reselect / re-activate the sheets the way the author left the workbook.
$wsAusgaben.Range("A4:I4").Select()
$ws.Activate()
$ws.Range("D11").Select()
